$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 115
$ws.Range("I2").Value = 306
$ws.Range("J2").Value = 1182
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 318
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = 209
$ws.Range("P2").Value = 9
$ws.Range("R2").Value = 12
$ws.Range("S2").Value = 128
$ws.Range("T2").Value = 213
$ws.Range("U2").Value = 17
$ws.Range("V2").Value = 1879
$ws.Range("X2").Value = 1910
$ws.Range("Y2").Value = 5
$ws.Range("Z2").Value = 25
$ws.Range("AA2").Value = 14
